$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 240, shifting existing rows 240-247 down to 241-248
$ws.Rows.Item(240).Insert()

# Populate the newly inserted row 240 with the new data
$ws.Cells.Item(240, 1).Value = 3
$ws.Cells.Item(240, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(240, 3).Value = "Coquimbo"
$ws.Cells.Item(240, 4).Value = 45075
$ws.Cells.Item(240, 5).Value = 5
$ws.Cells.Item(240, 6).Value = 100112026
$ws.Cells.Item(240, 7).Value = "Haba"
$ws.Cells.Item(240, 8).Value = "Sin especificar"
$ws.Cells.Item(240, 9).Value = "Primera"
$ws.Cells.Item(240, 10).Value = 85
$ws.Cells.Item(240, 11).Value = 19000
$ws.Cells.Item(240, 12).Value = 20000
$ws.Cells.Item(240, 13).Value = 19471
$ws.Cells.Item(240, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(240, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(240, 16).Value = 779
$ws.Cells.Item(240, 17).Value = 25
$ws.Cells.Item(240, 18).Value = "Hortaliza"
